# TC05_INS_CancerType-CervicalCancer.xlsx
# Update the "ProgramsTab" TabQuery (cell B2) so that the "Website" and
# "Data Location Details" columns are derived from prg.program_acronym
# (via a CASE on prg.program_link / prg.data_link) instead of prg.website.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newProgramsQuery = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.program_acronym     
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Cervical Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
'@

$ws.Range("B2").Value = $newProgramsQuery

# Re-apply the same visual formatting (wrapped, 12pt) so the cell keeps its
# existing look after the text swap.
$ws.Range("B2").Font.Name = "Calibri"
$ws.Range("B2").Font.Size = 12
$ws.Range("B2").WrapText = $true

# Move the saved selection/active cell from C5 to B2 and let the view
# scroll back so row 1 is visible again (drops the old topLeftCell="A5").
$ws.Range("B2").Select() | Out-Null
